$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2802.4443
$ws.Range("J17").Value = 2929.535
$ws.Range("L17").Value = 8788.605
$ws.Range("N17").Value = -9124.605
$ws.Range("H80").Value = 586.1579
$ws.Range("I80").Value = 574.2778
$ws.Range("J80").Value = 800
$ws.Range("K80").Value = 1722.8334
$ws.Range("L80").Value = 2400
$ws.Range("M80").Value = -724.8334
$ws.Range("N80").Value = -4396
$ws.Range("H83").Value = 586.1579
$ws.Range("I83").Value = 574.2778
$ws.Range("J83").Value = 800
$ws.Range("K83").Value = 5168.500199999999
$ws.Range("L83").Value = 7200
$ws.Range("M83").Value = -176.5001999999995
$ws.Range("N83").Value = -17184
$ws.Range("H137").Value = 2730.8572
$ws.Range("J137").Value = 1724.3334
$ws.Range("L137").Value = 5173.0002
$ws.Range("N137").Value = -10273.0002
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4387.2856
$ws.Range("I45").Value = 1553
$ws.Range("J45").Value = 8166.3335
$ws.Range("K45").Value = 1553
$ws.Range("L45").Value = 8166.3335
$ws.Range("M45").Value = -1176
$ws.Range("N45").Value = -8920.333500000001
$ws.Range("H63").Value = 4330.909
$ws.Range("I63").Value = 2380
$ws.Range("J63").Value = 5956.6665
$ws.Range("K63").Value = 2380
$ws.Range("L63").Value = 5956.6665
$ws.Range("M63").Value = -1694
$ws.Range("N63").Value = -7328.6665
$ws.Range("H66").Value = 4330.909
$ws.Range("I66").Value = 2380
$ws.Range("J66").Value = 5956.6665
$ws.Range("K66").Value = 11900
$ws.Range("L66").Value = 29783.3325
$ws.Range("M66").Value = -8468
$ws.Range("N66").Value = -36647.3325
$ws.Range("H88").Value = 2001.1666
$ws.Range("J88").Value = 2001.75
$ws.Range("L88").Value = 2001.75
$ws.Range("N88").Value = -2813.75
$ws.Range("H91").Value = 2001.1666
$ws.Range("J91").Value = 2001.75
$ws.Range("L91").Value = 2001.75
$ws.Range("N91").Value = -4809.75
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1378.3158
$ws.Range("I86").Value = 1399
$ws.Range("J86").Value = 1359.7
$ws.Range("K86").Value = 1399
$ws.Range("L86").Value = 1359.7
$ws.Range("M86").Value = -276
$ws.Range("N86").Value = -3605.7
$ws.Range("H89").Value = 1378.3158
$ws.Range("I89").Value = 1399
$ws.Range("J89").Value = 1359.7
$ws.Range("K89").Value = 6995
$ws.Range("L89").Value = 6798.5
$ws.Range("M89").Value = -1379
$ws.Range("N89").Value = -18030.5
$ws.Range("H105").Value = 3588.75
$ws.Range("I105").Value = 3932.5
$ws.Range("J105").Value = 3245
$ws.Range("K105").Value = 3932.5
$ws.Range("L105").Value = 3245
$ws.Range("M105").Value = -2185.5
$ws.Range("N105").Value = -6739
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9011338
$ws.Range("I31").Value = 1868
$ws.Range("K31").Value = 1868
$ws.Range("M31").Value = -1573
$ws.Range("H34").Value = 9011338
$ws.Range("I34").Value = 1868
$ws.Range("K34").Value = 1868
$ws.Range("M34").Value = -1666
$ws.Range("H58").Value = 795553
$ws.Range("I58").Value = 1232.449
$ws.Range("J58").Value = 3575675
$ws.Range("K58").Value = 1232.449
$ws.Range("L58").Value = 3575675
$ws.Range("M58").Value = -1029.449
$ws.Range("N58").Value = -3576081
$ws.Range("H122").Value = 76925950
$ws.Range("I122").Value = 90909760
$ws.Range("J122").Value = 15000
$ws.Range("K122").Value = 272729280
$ws.Range("L122").Value = 45000
$ws.Range("M122").Value = -272726830
$ws.Range("N122").Value = -49900
$ws.Range("H132").Value = 2907.3333
$ws.Range("I132").Value = 1877.8
$ws.Range("J132").Value = 3642.7144
$ws.Range("K132").Value = 5633.4
$ws.Range("L132").Value = 10928.1432
$ws.Range("M132").Value = -3103.4
$ws.Range("N132").Value = -15988.1432
$ws.Range("H136").Value = 795553
$ws.Range("I136").Value = 1232.449
$ws.Range("J136").Value = 3575675
$ws.Range("K136").Value = 3697.347
$ws.Range("L136").Value = 10727025
$ws.Range("M136").Value = -1147.347
$ws.Range("N136").Value = -10732125
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 6145.8623
$ws.Range("I125").Value = 3207.5
$ws.Range("J125").Value = 6616
$ws.Range("K125").Value = 9622.5
$ws.Range("L125").Value = 19848
$ws.Range("M125").Value = -4702.5
$ws.Range("N125").Value = -29688
$ws.Range("H132").Value = 3342.2917
$ws.Range("J132").Value = 8572.143
$ws.Range("L132").Value = 77149.287
$ws.Range("N132").Value = -82209.287
$ws.Range("H141").Value = 4443.684
$ws.Range("I141").Value = 3157.5
$ws.Range("J141").Value = 4786.6665
$ws.Range("K141").Value = 9472.5
$ws.Range("L141").Value = 14359.9995
$ws.Range("M141").Value = -4292.5
$ws.Range("N141").Value = -24719.9995
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 335788.94
$ws.Range("I122").Value = 386641.84
$ws.Range("J122").Value = 5245
$ws.Range("K122").Value = 1159925.52
$ws.Range("L122").Value = 15735
$ws.Range("M122").Value = -1157475.52
$ws.Range("N122").Value = -20635
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4331.75
$ws.Range("I61").Value = 1387
$ws.Range("J61").Value = 10221.25
$ws.Range("K61").Value = 1387
$ws.Range("L61").Value = 10221.25
$ws.Range("M61").Value = -1185
$ws.Range("N61").Value = -10625.25
$ws.Range("H93").Value = 1761.6774
$ws.Range("I93").Value = 1648.4445
$ws.Range("J93").Value = 2526
$ws.Range("K93").Value = 1648.4445
$ws.Range("L93").Value = 2526
$ws.Range("M93").Value = -400.4445000000001
$ws.Range("N93").Value = -5022
$ws.Range("H113").Value = 4331.75
$ws.Range("I113").Value = 1387
$ws.Range("J113").Value = 10221.25
$ws.Range("K113").Value = 1387
$ws.Range("L113").Value = 10221.25
$ws.Range("M113").Value = 783
$ws.Range("N113").Value = -14561.25
$ws.Range("H132").Value = 29414796
$ws.Range("I132").Value = 47621960
$ws.Range("J132").Value = 3222.923
$ws.Range("K132").Value = 142865880
$ws.Range("L132").Value = 9668.769
$ws.Range("M132").Value = -142863350
$ws.Range("N132").Value = -14728.769
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 21943.8
$ws.Range("I41").Value = 8342
$ws.Range("J41").Value = 25344.25
$ws.Range("K41").Value = 8342
$ws.Range("L41").Value = 25344.25
$ws.Range("M41").Value = -7952
$ws.Range("N41").Value = -26124.25
$ws.Range("H45").Value = 6158.625
$ws.Range("J45").Value = 6128.5713
$ws.Range("L45").Value = 6128.5713
$ws.Range("N45").Value = -7110.5713
$ws.Range("H74").Value = 6500
$ws.Range("I74").Value = 7000
$ws.Range("J74").Value = 6166.6665
$ws.Range("K74").Value = 7000
$ws.Range("L74").Value = 6166.6665
$ws.Range("M74").Value = -6064
$ws.Range("N74").Value = -8038.6665
$ws.Range("H77").Value = 6500
$ws.Range("I77").Value = 7000
$ws.Range("J77").Value = 6166.6665
$ws.Range("K77").Value = 21000
$ws.Range("L77").Value = 18499.9995
$ws.Range("M77").Value = -16320
$ws.Range("N77").Value = -27859.9995
$ws.Range("H96").Value = 8245.083000000001
$ws.Range("I96").Value = 3400.1667
$ws.Range("J96").Value = 13090
$ws.Range("K96").Value = 3400.1667
$ws.Range("L96").Value = 13090
$ws.Range("M96").Value = -2027.1667
$ws.Range("N96").Value = -15836
$ws.Range("H100").Value = 566.9167
$ws.Range("I100").Value = 510.94446
$ws.Range("J100").Value = 734.8333
$ws.Range("K100").Value = 1021.88892
$ws.Range("L100").Value = 1469.6666
$ws.Range("M100").Value = -480.88892
$ws.Range("N100").Value = -2551.6666
$ws.Range("H107").Value = 363.4
$ws.Range("I107").Value = 381.66666
$ws.Range("J107").Value = 199
$ws.Range("K107").Value = 1144.99998
$ws.Range("L107").Value = 597
$ws.Range("M107").Value = 775.0000199999999
$ws.Range("N107").Value = -4437
$ws.Range("H132").Value = 2916.375
$ws.Range("I132").Value = 2686.25
$ws.Range("J132").Value = 3031.4375
$ws.Range("K132").Value = 8058.75
$ws.Range("L132").Value = 9094.3125
$ws.Range("M132").Value = -5528.75
$ws.Range("N132").Value = -14154.3125
$ws.Range("H136").Value = 8335384
$ws.Range("I136").Value = 22728976
$ws.Range("J136").Value = 2251.7368
$ws.Range("K136").Value = 68186928
$ws.Range("L136").Value = 6755.2104
$ws.Range("M136").Value = -68184378
$ws.Range("N136").Value = -11855.2104
